$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.5871232612860219
$ws.Range("C2").Value = -0.7413052565746122
$ws.Range("D2").Value = -0.6779821139628318

$ws.Range("B3").Value = 0.6848673230507893
$ws.Range("C3").Value = 0.6903358201961617
$ws.Range("D3").Value = 0.6978903635809913

$ws.Range("B4").Value = 0.6476181060054425
$ws.Range("C4").Value = -0.7068555247253859
$ws.Range("D4").Value = 0.7256394420805552

$ws.Range("B5").Value = -0.8036534091992149
$ws.Range("C5").Value = -0.5697401020314303
$ws.Range("D5").Value = 0.7061685519639862

$ws.Range("B6").Value = 0.5138994230556155
$ws.Range("C6").Value = -0.6515289144268095
$ws.Range("D6").Value = -0.6137815080324943

$ws.Range("B7").Value = -0.7696509785612776
$ws.Range("C7").Value = 0.6610223880034173
$ws.Range("D7").Value = 0.7415631069330682

$ws.Range("B8").Value = 0.6805028901862737
$ws.Range("C8").Value = -0.6126335702838159
$ws.Range("D8").Value = -0.618224751893086

$ws.Range("B9").Value = 0.6972657875223104
$ws.Range("C9").Value = 0.7828012631545163
$ws.Range("D9").Value = -0.6113732004916022
